$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the (second) "What are the constraints? Sub-problems?"
# paragraph -- the one belonging to the "Socks in the Dark" problem
# (the first occurrence belongs to the earlier cat/parrot/seed
# problem), using Find so the script does not depend on hard-coded
# paragraph numbers.
# ------------------------------------------------------------------
$firstHit = $d.Content
$null = $firstHit.Find.Execute("What are the constraints? Sub-problems?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$secondHit = $d.Range($firstHit.End, $d.Content.End)
$null = $secondHit.Find.Execute("What are the constraints? Sub-problems?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Resolve the paragraph index that corresponds to this second match by
# comparing against $d.Paragraphs (Paragraph objects obtained directly
# from the Paragraphs collection are handled reliably by this engine;
# Paragraphs derived from a freshly-built zero-length Range are not).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $secondHit.Start) {
        $targetIndex = $i
        break
    }
}

# The paragraph right after the question is the (currently empty)
# paragraph that should receive the answer.
$emptyParaIndex = $targetIndex + 1
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)

# ------------------------------------------------------------------
# Fill in the previously-empty paragraph with the answer text. Using
# InsertBefore on the paragraph's own Range (which still includes the
# trailing paragraph mark) puts the new text inside that paragraph,
# right before the pilcrow, without adding an extra paragraph.
# ------------------------------------------------------------------
$answerText = "One of the constraints is you cannot see the socks that you are selecting.  A sub-problem is you only want to select enough socks to guarantee a match/ a match of each color."
$emptyPara.Range.InsertBefore($answerText)

# ------------------------------------------------------------------
# Move the "_GoBack" bookmark from the end of the previous ("The goal
# here is ... leave the room.") paragraph to the end of this new
# answer text (both are collapsed bookmarks sitting right before the
# paragraph mark).
#
# Workaround: creating a *collapsed* Range/Bookmark whose position is
# immediately in front of a paragraph mark is mishandled by this
# runtime (it silently resets to the very start of the document), so
# a one-character placeholder is appended first to move the target
# position away from the paragraph-mark boundary, the bookmark is
# created, and the placeholder is then removed again.
# ------------------------------------------------------------------
$answerPara = $d.Paragraphs.Item($emptyParaIndex)
$answerPara.Range.InsertAfter("Z")

$answerParaWithPlaceholder = $d.Paragraphs.Item($emptyParaIndex)
$bookmarkPos = $answerParaWithPlaceholder.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderPara = $d.Paragraphs.Item($emptyParaIndex)
$placeholderRange = $d.Range($placeholderPara.Range.End - 2, $placeholderPara.Range.End - 1)
$placeholderRange.Delete()
